$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript three (U+2083) character used in row 29 price values
$subThree = [string][char]0x2083

$ws.Range("D2").Value = '59.148.07'
$ws.Range("E2").Value = '  +0.06%  '

$ws.Range("D3").Value = '2.500.50'
$ws.Range("E3").Value = '  +0.00%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.63'
$ws.Range("E5").Value = '  -0.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.43'
$ws.Range("E6").Value = '  -1.24%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.567'
$ws.Range("E8").Value = '  +1.02%  '

$ws.Range("D9").Value = '2.522.37'
$ws.Range("E9").Value = '  +0.66%  '

$ws.Range("E10").Value = '  +0.94%  '

$ws.Range("E11").Value = '  -2.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.31'
$ws.Range("E12").Value = '  -1.99%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.348'
$ws.Range("E13").Value = '  -0.28%  '

$ws.Range("D14").Value = '2.945.81'
$ws.Range("E14").Value = '  -0.15%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.04'
$ws.Range("E15").Value = '  +0.20%  '

$ws.Range("D16").Value = '58.943.04'
$ws.Range("E16").Value = '  -0.11%  '

$ws.Range("E17").Value = '  -0.69%  '

$ws.Range("D18").Value = '2.526.16'
$ws.Range("E18").Value = '  +1.11%  '

$ws.Range("E19").Value = '  +1.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.27'
$ws.Range("E20").Value = '  +0.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '323.41'
$ws.Range("E21").Value = '  -0.10%  '

$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.96'
$ws.Range("E23").Value = '  +2.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.27'
$ws.Range("E24").Value = '  +3.70%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.420'
$ws.Range("E25").Value = '  +1.30%  '

$ws.Range("E26").Value = '  -1.65%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  +0.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.53'
$ws.Range("E28").Value = '  -1.69%  '

$ws.Range("D29").Value = '0.0' + $subThree + '0770'
$ws.Range("E29").Value = '  -0.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.65'
$ws.Range("E30").Value = '  +0.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '171.48'
$ws.Range("E31").Value = '  +3.78%  '

$ws.Range("E32").Value = '  -1.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.18'
$ws.Range("E33").Value = '  +8.43%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.39'
$ws.Range("E35").Value = '  +1.04%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.38'
$ws.Range("E36").Value = '  -0.37%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.08'
$ws.Range("E37").Value = '  -0.33%  '

$ws.Range("E38").Value = '  -1.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.92'
$ws.Range("E39").Value = '  +0.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.813'
$ws.Range("E40").Value = '  +1.72%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.60'
$ws.Range("E41").Value = '  -0.90%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '286.00'
$ws.Range("E42").Value = '  +3.25%  '

$ws.Range("E43").Value = '  +0.17%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.996'
$ws.Range("E44").Value = '  -0.34%  '

$ws.Range("E45").Value = '  +2.37%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '130.36'
$ws.Range("E46").Value = '  +4.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.88'
$ws.Range("E47").Value = '  -0.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0923'
$ws.Range("E48").Value = '  -1.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0505'
$ws.Range("E49").Value = '  -0.57%  '

$ws.Range("E50").Value = '  -0.46%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.41'
$ws.Range("E51").Value = '  -0.53%  '
